$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 (CAMERA_DISTANCE_STRENGTH / 20.f) is unchanged - leave as-is.

# Update row numbers / labels for the existing camera-distance rows and
# insert two brand new rows (Market min/max), in the precise order the
# strings first appear so the shared-string table matches the target file.

# Row 17: MAX_CAMERA_DISTANCE -> MAX_CAMERA_DISTANCE_TOWN
$ws.Cells.Item(17,1).Value = 15
$ws.Cells.Item(17,2).Value = "MAX_CAMERA_DISTANCE_TOWN"

# Row 18: MIN_CAMERA_DISTANCE -> MIN_CAMERA_DISTANCE_TOWN, clear stray C18 value
$ws.Cells.Item(18,1).Value = 16
$ws.Cells.Item(18,2).Value = "MIN_CAMERA_DISTANCE_TOWN"
$ws.Cells.Item(18,3).ClearContents()
$ws.Cells.Item(18,4).Value = "300.f"

# Row 17 value: 2000.f -> 1000.f
$ws.Cells.Item(17,4).Value = "1000.f"

# New row 19 (Market max) - set D first so "1500.f" registers before the
# MAX_CAMERA_DISTANCE_MARKET label, matching the authored string order.
$ws.Cells.Item(19,1).Value = 17
$ws.Cells.Item(19,4).Value = "1500.f"
$ws.Cells.Item(19,2).Value = "MAX_CAMERA_DISTANCE_MARKET"

# New row 20 (Market min)
$ws.Cells.Item(20,1).Value = 18
$ws.Cells.Item(20,2).Value = "MIN_CAMERA_DISTANCE_MARKET"
$ws.Cells.Item(20,4).Value = "300.f"

# Match styling of rows above: A and B columns use the vertical-center
# alignment style (same as rows 16-18); C is left untouched (no cell at
# all) and D keeps the default style, mirroring the shape of rows 16-18.
$ws.Range("A19:B20").VerticalAlignment = -4108

# Update the active selection to match the saved file.
$ws.Range("A17:D20").Select()
